# "added initial file input design to admin page"
#
# Target slide is the Admin page (p:sldId cId/sldId 259 / creationId
# 2097583108), which is slide index 5 in the presentation's slide order.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(5)

# ---------------------------------------------------------------------
# 1. Nudge the full-bleed background rectangle (id 4, "Rectangle 3").
# ---------------------------------------------------------------------
$rect3 = $s.Shapes.Item(1)          # id 4 "Rectangle 3"
$rect3.Left = 9.171969503937007
$rect3.Top  = 1.829685039370079

# ---------------------------------------------------------------------
# 2. Nudge the main content panel (id 3, "Rectangle 2").
# ---------------------------------------------------------------------
$rect2 = $s.Shapes.Item(7)          # id 3 "Rectangle 2"
$rect2.Left = 41.426851393700787
$rect2.Top  = 96.926456692913391

# ---------------------------------------------------------------------
# 3. Move the existing "Upload match CSV" button up and rename it to
#    "Upload Season" (id 11, "Rounded Rectangle 10").
# ---------------------------------------------------------------------
$uploadSeason = $s.Shapes.Item(8)   # id 11 "Rounded Rectangle 10"
$uploadSeason.Left = 309.060472440944864
$uploadSeason.Top  = 112.756299212598421
$uploadSeason.TextFrame.TextRange.Text = "Upload Season"

# ---------------------------------------------------------------------
# 4. Shrink/move the Message Box panel (id 13, "Rectangle 12").
# ---------------------------------------------------------------------
$msgBox = $s.Shapes.Item(9)         # id 13 "Rectangle 12"
$msgBox.Left   = 736.897979755905567
$msgBox.Top    = 96.926456692913391
$msgBox.Width  = 181.675196850393689
$msgBox.Height = 289.780393700787386

# ---------------------------------------------------------------------
# 5. Add the new admin-page buttons by duplicating the "Upload Season"
#    rounded-rectangle button (keeps its theme style/geometry/name) and
#    repositioning + relabelling each copy. Duplicating in this order
#    reproduces the shape ids from the authored edit (2, 10, 12, 14, 15).
# ---------------------------------------------------------------------

# id 2 -> "Add new user"
$addUser = $uploadSeason.Duplicate()
$addUser.Left   = 736.703779527559050
$addUser.Top    = 406.588669417322819
$addUser.Width  = 181.675196850393689
$addUser.Height = 40.976614173228349
$addUser.TextFrame.TextRange.Text = "Add new user"

# id 10 -> "Upload Players"
$uploadPlayers = $uploadSeason.Duplicate()
$uploadPlayers.Left   = 309.060472440944864
$uploadPlayers.Top    = 175.166456692913385
$uploadPlayers.Width  = 127.948979377952767
$uploadPlayers.Height = 43.108267716535430
$uploadPlayers.TextFrame.TextRange.Text = "Upload Players"

# id 12 -> "Upload Prize Money"
$uploadPrize = $uploadSeason.Duplicate()
$uploadPrize.Left   = 309.060472440944864
$uploadPrize.Top    = 241.816614173228345
$uploadPrize.Width  = 127.948979377952767
$uploadPrize.Height = 43.108267716535430
$uploadPrize.TextFrame.TextRange.Text = "Upload Prize Money"

# id 14 -> "Submit"
$submit = $uploadSeason.Duplicate()
$submit.Left   = 309.060472440944864
$submit.Top    = 375.116929133858264
$submit.Width  = 127.948979377952767
$submit.Height = 43.108267716535430
$submit.TextFrame.TextRange.Text = "Submit"

# id 15 -> "Upload Matches"
$uploadMatches = $uploadSeason.Duplicate()
$uploadMatches.Left   = 309.060472440944864
$uploadMatches.Top    = 308.466781653543308
$uploadMatches.Width  = 127.948979377952767
$uploadMatches.Height = 43.108267716535430
$uploadMatches.TextFrame.TextRange.Text = "Upload Matches"
